$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 1690.4517
$ws.Range("I40").Value = 1550.4286
$ws.Range("K40").Value = 1550.4286
$ws.Range("M40").Value = -1375.4286
$ws.Range("H63").Value = 46666.668
$ws.Range("I63").Value = 20000
$ws.Range("J63").Value = 60000
$ws.Range("K63").Value = 20000
$ws.Range("L63").Value = 60000
$ws.Range("M63").Value = -19376
$ws.Range("N63").Value = -61248
$ws.Range("H66").Value = 46666.668
$ws.Range("I66").Value = 20000
$ws.Range("J66").Value = 60000
$ws.Range("K66").Value = 60000
$ws.Range("L66").Value = 180000
$ws.Range("M66").Value = -56880
$ws.Range("N66").Value = -186240
$ws.Range("H100").Value = 1432.6111
$ws.Range("I100").Value = 1399.2354
$ws.Range("K100").Value = 1399.2354
$ws.Range("M100").Value = -858.2354
$ws.Range("H107").Value = 654.0909
$ws.Range("I107").Value = 667.1905
$ws.Range("J107").Value = 379
$ws.Range("K107").Value = 667.1905
$ws.Range("L107").Value = 379
$ws.Range("M107").Value = 1252.8095
$ws.Range("N107").Value = -4219
$ws.Range("H112").Value = 75111.07000000001
$ws.Range("J112").Value = 80860
$ws.Range("L112").Value = 242580
$ws.Range("N112").Value = -244796
$ws.Range("H138").Value = 4020.4211
$ws.Range("J138").Value = 5098.385
$ws.Range("L138").Value = 15295.155
$ws.Range("N138").Value = -25575.155
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 5607.0303
$ws.Range("I45").Value = 9252.857
$ws.Range("K45").Value = 9252.857
$ws.Range("M45").Value = -8875.857
$ws.Range("H61").Value = 4051.6316
$ws.Range("I61").Value = 3608.0625
$ws.Range("K61").Value = 3608.0625
$ws.Range("M61").Value = -3396.0625
$ws.Range("H102").Value = 4381.8667
$ws.Range("I102").Value = 4473.25
$ws.Range("K102").Value = 4473.25
$ws.Range("M102").Value = -2851.25
$ws.Range("H132").Value = 3218.16
$ws.Range("I132").Value = 2668.55
$ws.Range("J132").Value = 5416.6
$ws.Range("K132").Value = 8005.650000000001
$ws.Range("L132").Value = 16249.8
$ws.Range("M132").Value = -5475.650000000001
$ws.Range("N132").Value = -21309.8
$ws.Range("H136").Value = 4051.6316
$ws.Range("I136").Value = 3608.0625
$ws.Range("K136").Value = 10824.1875
$ws.Range("M136").Value = -8274.1875
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 356430.06
$ws.Range("I22").Value = 512.5833
$ws.Range("K22").Value = 512.5833
$ws.Range("M22").Value = -339.5833
$ws.Range("H105").Value = 3594.6667
$ws.Range("J105").Value = 2300
$ws.Range("L105").Value = 2300
$ws.Range("N105").Value = -5794
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 38107.45
$ws.Range("I31").Value = 49003.855
$ws.Range("K31").Value = 49003.855
$ws.Range("M31").Value = -48708.855
$ws.Range("H34").Value = 38107.45
$ws.Range("I34").Value = 49003.855
$ws.Range("K34").Value = 49003.855
$ws.Range("M34").Value = -48801.855
$ws.Range("H58").Value = 2742.647
$ws.Range("J58").Value = 2095
$ws.Range("L58").Value = 2095
$ws.Range("N58").Value = -2501
$ws.Range("H63").Value = 0
$ws.Range("J63").Value = 0
$ws.Range("L63").Value = 0
$ws.Range("N63").ClearContents()
$ws.Range("H64").Value = 0
$ws.Range("J64").Value = 0
$ws.Range("L64").Value = 0
$ws.Range("N64").ClearContents()
$ws.Range("H66").Value = 0
$ws.Range("J66").Value = 0
$ws.Range("L66").Value = 0
$ws.Range("N66").ClearContents()
$ws.Range("H67").Value = 0
$ws.Range("J67").Value = 0
$ws.Range("L67").Value = 0
$ws.Range("N67").ClearContents()
$ws.Range("H99").Value = 2839.9
$ws.Range("I99").Value = 2839.9
$ws.Range("K99").Value = 2839.9
$ws.Range("M99").Value = -1341.9
$ws.Range("H126").Value = 2839.9
$ws.Range("I126").Value = 2839.9
$ws.Range("K126").Value = 8519.700000000001
$ws.Range("M126").Value = -6049.700000000001
$ws.Range("H132").Value = 3229.6667
$ws.Range("I132").Value = 2883.375
$ws.Range("K132").Value = 8650.125
$ws.Range("M132").Value = -6120.125
$ws.Range("H134").Value = 7626.7856
$ws.Range("I134").Value = 6249.548
$ws.Range("K134").Value = 18748.644
$ws.Range("M134").Value = -16213.644
$ws.Range("H136").Value = 2742.647
$ws.Range("J136").Value = 2095
$ws.Range("L136").Value = 6285
$ws.Range("N136").Value = -11385
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H124").Value = 28463.264
$ws.Range("I124").Value = 6956
$ws.Range("J124").Value = 32495.875
$ws.Range("K124").Value = 20868
$ws.Range("L124").Value = 97487.625
$ws.Range("M124").Value = -15958
$ws.Range("N124").Value = -107307.625
$ws.Range("H128").Value = 349861.44
$ws.Range("I128").Value = 349861.44
$ws.Range("K128").Value = 1049584.32
$ws.Range("M128").Value = -1044604.32
$ws.Range("H129").Value = 756.2
$ws.Range("I129").Value = 662
$ws.Range("J129").Value = 1133
$ws.Range("K129").Value = 1986
$ws.Range("L129").Value = 3399
$ws.Range("M129").Value = 3014
$ws.Range("N129").Value = -13399
$ws.Range("H131").Value = 56946.668
$ws.Range("J131").Value = 2491.4285
$ws.Range("L131").Value = 7474.2855
$ws.Range("N131").Value = -17554.2855
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H64").Value = 32000
$ws.Range("J64").Value = 32000
$ws.Range("L64").Value = 32000
$ws.Range("N64").Value = -32496
$ws.Range("H67").Value = 32000
$ws.Range("J67").Value = 32000
$ws.Range("L67").Value = 32000
$ws.Range("N67").Value = -33716
$ws.Range("H102").Value = 2138.2727
$ws.Range("I102").Value = 2151.4
$ws.Range("K102").Value = 2151.4
$ws.Range("M102").Value = -529.4000000000001
$ws.Range("H132").Value = 3283.2856
$ws.Range("I132").Value = 3097.45
$ws.Range("K132").Value = 9292.349999999999
$ws.Range("M132").Value = -6762.349999999999
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1352.1
$ws.Range("I22").Value = 1310
$ws.Range("J22").Value = 1370.1428
$ws.Range("K22").Value = 1310
$ws.Range("L22").Value = 1370.1428
$ws.Range("M22").Value = -1015
$ws.Range("N22").Value = -1960.1428
$ws.Range("H27").Value = 1352.1
$ws.Range("I27").Value = 1310
$ws.Range("J27").Value = 1370.1428
$ws.Range("K27").Value = 1310
$ws.Range("L27").Value = 1370.1428
$ws.Range("M27").Value = -1203
$ws.Range("N27").Value = -1584.1428
$ws.Range("H82").Value = 2172
$ws.Range("I82").Value = 2172
$ws.Range("K82").Value = 2172
$ws.Range("M82").Value = -1811
$ws.Range("H85").Value = 2172
$ws.Range("I85").Value = 2172
$ws.Range("K85").Value = 2172
$ws.Range("M85").Value = -924
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 3080.162
$ws.Range("I132").Value = 3210.9656
$ws.Range("J132").Value = 2606
$ws.Range("K132").Value = 9632.8968
$ws.Range("L132").Value = 7818
$ws.Range("M132").Value = -7102.8968
$ws.Range("N132").Value = -12878
